$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-21 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("29×14=", $true, $false, $false, $false, $false, $true, 1, $false, "48×80=", 2) | Out-Null
$d.Content.Find.Execute("60×19=", $true, $false, $false, $false, $false, $true, 1, $false, "99×26=", 2) | Out-Null
$d.Content.Find.Execute("49×47=", $true, $false, $false, $false, $false, $true, 1, $false, "69×99=", 2) | Out-Null
$d.Content.Find.Execute("91×92=", $true, $false, $false, $false, $false, $true, 1, $false, "57×38=", 2) | Out-Null
$d.Content.Find.Execute("84×23=", $true, $false, $false, $false, $false, $true, 1, $false, "22×63=", 2) | Out-Null
$d.Content.Find.Execute("33×12=", $true, $false, $false, $false, $false, $true, 1, $false, "72×14=", 2) | Out-Null
$d.Content.Find.Execute("65×99=", $true, $false, $false, $false, $false, $true, 1, $false, "84×98=", 2) | Out-Null
$d.Content.Find.Execute("96×18=", $true, $false, $false, $false, $false, $true, 1, $false, "40×96=", 2) | Out-Null
$d.Content.Find.Execute("13×94=", $true, $false, $false, $false, $false, $true, 1, $false, "13×91=", 2) | Out-Null
$d.Content.Find.Execute("78×33=", $true, $false, $false, $false, $false, $true, 1, $false, "35×99=", 2) | Out-Null
$d.Content.Find.Execute("64×81=", $true, $false, $false, $false, $false, $true, 1, $false, "62×90=", 2) | Out-Null
$d.Content.Find.Execute("50×75=", $true, $false, $false, $false, $false, $true, 1, $false, "46×70=", 2) | Out-Null
$d.Content.Find.Execute("79×14=", $true, $false, $false, $false, $false, $true, 1, $false, "57×43=", 2) | Out-Null
$d.Content.Find.Execute("85×38=", $true, $false, $false, $false, $false, $true, 1, $false, "42×20=", 2) | Out-Null
$d.Content.Find.Execute("42×38=", $true, $false, $false, $false, $false, $true, 1, $false, "72×92=", 2) | Out-Null
$d.Content.Find.Execute("21×39=", $true, $false, $false, $false, $false, $true, 1, $false, "95×98=", 2) | Out-Null
$d.Content.Find.Execute("46×39=", $true, $false, $false, $false, $false, $true, 1, $false, "40×55=", 2) | Out-Null
$d.Content.Find.Execute("14×90=", $true, $false, $false, $false, $false, $true, 1, $false, "90×58=", 2) | Out-Null
$d.Content.Find.Execute("96×52=", $true, $false, $false, $false, $false, $true, 1, $false, "63×49=", 2) | Out-Null
$d.Content.Find.Execute("29×10=", $true, $false, $false, $false, $false, $true, 1, $false, "73×37=", 2) | Out-Null
$d.Content.Find.Execute("31×98=", $true, $false, $false, $false, $false, $true, 1, $false, "76×86=", 2) | Out-Null
$d.Content.Find.Execute("43×92=", $true, $false, $false, $false, $false, $true, 1, $false, "73×73=", 2) | Out-Null
$d.Content.Find.Execute("16×19=", $true, $false, $false, $false, $false, $true, 1, $false, "34×55=", 2) | Out-Null
$d.Content.Find.Execute("66×85=", $true, $false, $false, $false, $false, $true, 1, $false, "79×86=", 2) | Out-Null
$d.Content.Find.Execute("82×18=", $true, $false, $false, $false, $false, $true, 1, $false, "30×67=", 2) | Out-Null
$d.Content.Find.Execute("75×37=", $true, $false, $false, $false, $false, $true, 1, $false, "88×56=", 2) | Out-Null
$d.Content.Find.Execute("28×55=", $true, $false, $false, $false, $false, $true, 1, $false, "74×90=", 2) | Out-Null
$d.Content.Find.Execute("71×29=", $true, $false, $false, $false, $false, $true, 1, $false, "46×96=", 2) | Out-Null
$d.Content.Find.Execute("36×21=", $true, $false, $false, $false, $false, $true, 1, $false, "87×100=", 2) | Out-Null
$d.Content.Find.Execute("93×97=", $true, $false, $false, $false, $false, $true, 1, $false, "45×34=", 2) | Out-Null
$d.Content.Find.Execute("15×58=", $true, $false, $false, $false, $false, $true, 1, $false, "36×42=", 2) | Out-Null
$d.Content.Find.Execute("67×79=", $true, $false, $false, $false, $false, $true, 1, $false, "69×19=", 2) | Out-Null
$d.Content.Find.Execute("24×40=", $true, $false, $false, $false, $false, $true, 1, $false, "97×10=", 2) | Out-Null
$d.Content.Find.Execute("35×12=", $true, $false, $false, $false, $false, $true, 1, $false, "69×25=", 2) | Out-Null
$d.Content.Find.Execute("14×95=", $true, $false, $false, $false, $false, $true, 1, $false, "39×28=", 2) | Out-Null
$d.Content.Find.Execute("55×14=", $true, $false, $false, $false, $false, $true, 1, $false, "47×77=", 2) | Out-Null
$d.Content.Find.Execute("62×29=", $true, $false, $false, $false, $false, $true, 1, $false, "22×32=", 2) | Out-Null
$d.Content.Find.Execute("16×35=", $true, $false, $false, $false, $false, $true, 1, $false, "51×41=", 2) | Out-Null
$d.Content.Find.Execute("26×40=", $true, $false, $false, $false, $false, $true, 1, $false, "26×10=", 2) | Out-Null
$d.Content.Find.Execute("82×84=", $true, $false, $false, $false, $false, $true, 1, $false, "88×17=", 2) | Out-Null
$d.Content.Find.Execute("29×97=", $true, $false, $false, $false, $false, $true, 1, $false, "93×78=", 2) | Out-Null
$d.Content.Find.Execute("62×51=", $true, $false, $false, $false, $false, $true, 1, $false, "40×45=", 2) | Out-Null
$d.Content.Find.Execute("97×50=", $true, $false, $false, $false, $false, $true, 1, $false, "30×51=", 2) | Out-Null
$d.Content.Find.Execute("53×58=", $true, $false, $false, $false, $false, $true, 1, $false, "58×72=", 2) | Out-Null
$d.Content.Find.Execute("44×71=", $true, $false, $false, $false, $false, $true, 1, $false, "90×60=", 2) | Out-Null
$d.Content.Find.Execute("70×11=", $true, $false, $false, $false, $false, $true, 1, $false, "67×55=", 2) | Out-Null
$d.Content.Find.Execute("97×62=", $true, $false, $false, $false, $false, $true, 1, $false, "12×47=", 2) | Out-Null
$d.Content.Find.Execute("72×85=", $true, $false, $false, $false, $false, $true, 1, $false, "25×98=", 2) | Out-Null
$d.Content.Find.Execute("97×91=", $true, $false, $false, $false, $false, $true, 1, $false, "20×18=", 2) | Out-Null
$d.Content.Find.Execute("98×82=", $true, $false, $false, $false, $false, $true, 1, $false, "72×82=", 2) | Out-Null
$d.Content.Find.Execute("50×52=", $true, $false, $false, $false, $false, $true, 1, $false, "78×50=", 2) | Out-Null
$d.Content.Find.Execute("70×53=", $true, $false, $false, $false, $false, $true, 1, $false, "82×19=", 2) | Out-Null
$d.Content.Find.Execute("13×42=", $true, $false, $false, $false, $false, $true, 1, $false, "66×44=", 2) | Out-Null
$d.Content.Find.Execute("97×25=", $true, $false, $false, $false, $false, $true, 1, $false, "73×61=", 2) | Out-Null
$d.Content.Find.Execute("96×85=", $true, $false, $false, $false, $false, $true, 1, $false, "26×58=", 2) | Out-Null
$d.Content.Find.Execute("42×23=", $true, $false, $false, $false, $false, $true, 1, $false, "33×77=", 2) | Out-Null
$d.Content.Find.Execute("100×21=", $true, $false, $false, $false, $false, $true, 1, $false, "20×40=", 2) | Out-Null
$d.Content.Find.Execute("85×39=", $true, $false, $false, $false, $false, $true, 1, $false, "71×61=", 2) | Out-Null
$d.Content.Find.Execute("30×80=", $true, $false, $false, $false, $false, $true, 1, $false, "25×14=", 2) | Out-Null
$d.Content.Find.Execute("85×16=", $true, $false, $false, $false, $false, $true, 1, $false, "52×54=", 2) | Out-Null
$d.Content.Find.Execute("69×47=", $true, $false, $false, $false, $false, $true, 1, $false, "47×17=", 2) | Out-Null
$d.Content.Find.Execute("71×30=", $true, $false, $false, $false, $false, $true, 1, $false, "43×62=", 2) | Out-Null
$d.Content.Find.Execute("65×38=", $true, $false, $false, $false, $false, $true, 1, $false, "79×20=", 2) | Out-Null
$d.Content.Find.Execute("80×71=", $true, $false, $false, $false, $false, $true, 1, $false, "26×69=", 2) | Out-Null
$d.Content.Find.Execute("93×41=", $true, $false, $false, $false, $false, $true, 1, $false, "83×87=", 2) | Out-Null
$d.Content.Find.Execute("18×22=", $true, $false, $false, $false, $false, $true, 1, $false, "69×57=", 2) | Out-Null
$d.Content.Find.Execute("98×45=", $true, $false, $false, $false, $false, $true, 1, $false, "71×79=", 2) | Out-Null
$d.Content.Find.Execute("65×94=", $true, $false, $false, $false, $false, $true, 1, $false, "79×21=", 2) | Out-Null
$d.Content.Find.Execute("17×21=", $true, $false, $false, $false, $false, $true, 1, $false, "80×64=", 2) | Out-Null
$d.Content.Find.Execute("58×67=", $true, $false, $false, $false, $false, $true, 1, $false, "59×76=", 2) | Out-Null
$d.Content.Find.Execute("24×92=", $true, $false, $false, $false, $false, $true, 1, $false, "61×72=", 2) | Out-Null
$d.Content.Find.Execute("27×17=", $true, $false, $false, $false, $false, $true, 1, $false, "43×89=", 2) | Out-Null
$d.Content.Find.Execute("27×35=", $true, $false, $false, $false, $false, $true, 1, $false, "18×92=", 2) | Out-Null
$d.Content.Find.Execute("82×93=", $true, $false, $false, $false, $false, $true, 1, $false, "82×80=", 2) | Out-Null
$d.Content.Find.Execute("72×29=", $true, $false, $false, $false, $false, $true, 1, $false, "49×34=", 2) | Out-Null
$d.Content.Find.Execute("36×41=", $true, $false, $false, $false, $false, $true, 1, $false, "30×30=", 2) | Out-Null
$d.Content.Find.Execute("72×33=", $true, $false, $false, $false, $false, $true, 1, $false, "68×22=", 2) | Out-Null
$d.Content.Find.Execute("46×88=", $true, $false, $false, $false, $false, $true, 1, $false, "31×68=", 2) | Out-Null
$d.Content.Find.Execute("67×22=", $true, $false, $false, $false, $false, $true, 1, $false, "83×22=", 2) | Out-Null
$d.Content.Find.Execute("98×64=", $true, $false, $false, $false, $false, $true, 1, $false, "94×78=", 2) | Out-Null
$d.Content.Find.Execute("85×65=", $true, $false, $false, $false, $false, $true, 1, $false, "18×51=", 2) | Out-Null
$d.Content.Find.Execute("83×15=", $true, $false, $false, $false, $false, $true, 1, $false, "70×12=", 2) | Out-Null
$d.Content.Find.Execute("90×21=", $true, $false, $false, $false, $false, $true, 1, $false, "57×68=", 2) | Out-Null
$d.Content.Find.Execute("100×58=", $true, $false, $false, $false, $false, $true, 1, $false, "61×20=", 2) | Out-Null
$d.Content.Find.Execute("46×63=", $true, $false, $false, $false, $false, $true, 1, $false, "70×28=", 2) | Out-Null
$d.Content.Find.Execute("45×69=", $true, $false, $false, $false, $false, $true, 1, $false, "45×65=", 2) | Out-Null
$d.Content.Find.Execute("13×13=", $true, $false, $false, $false, $false, $true, 1, $false, "54×68=", 2) | Out-Null
$d.Content.Find.Execute("43×64=", $true, $false, $false, $false, $false, $true, 1, $false, "83×36=", 2) | Out-Null
$d.Content.Find.Execute("90×93=", $true, $false, $false, $false, $false, $true, 1, $false, "19×39=", 2) | Out-Null
$d.Content.Find.Execute("10×67=", $true, $false, $false, $false, $false, $true, 1, $false, "89×68=", 2) | Out-Null
$d.Content.Find.Execute("56×61=", $true, $false, $false, $false, $false, $true, 1, $false, "26×46=", 2) | Out-Null
$d.Content.Find.Execute("78×91=", $true, $false, $false, $false, $false, $true, 1, $false, "56×60=", 2) | Out-Null
$d.Content.Find.Execute("63×81=", $true, $false, $false, $false, $false, $true, 1, $false, "56×89=", 2) | Out-Null
$d.Content.Find.Execute("44×38=", $true, $false, $false, $false, $false, $true, 1, $false, "16×56=", 2) | Out-Null
$d.Content.Find.Execute("32×68=", $true, $false, $false, $false, $false, $true, 1, $false, "52×55=", 2) | Out-Null
$d.Content.Find.Execute("13×81=", $true, $false, $false, $false, $false, $true, 1, $false, "63×52=", 2) | Out-Null
$d.Content.Find.Execute("20×26=", $true, $false, $false, $false, $false, $true, 1, $false, "42×85=", 2) | Out-Null
$d.Content.Find.Execute("61×30=", $true, $false, $false, $false, $false, $true, 1, $false, "80×45=", 2) | Out-Null
$d.Content.Find.Execute("88×64=", $true, $false, $false, $false, $false, $true, 1, $false, "79×84=", 2) | Out-Null
$d.Content.Find.Execute("46×29=", $true, $false, $false, $false, $false, $true, 1, $false, "73×12=", 2) | Out-Null
